$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update rule row C10 ("C1" condition threshold for rule R30) from 18 to 100
$ws.Range("C10").Value = 100
